$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3 was text "c"; change it to a numeric value (general numeric format handling)
$ws.Range("B3").Value = 2.66666156237642

# Move the active selection to B3
$ws.Range("B3").Select()
